$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run results for the existing "Full" model changed slightly too (pD/DIC updated).
$ws.Range("D2").Value2 = 80.825
$ws.Range("E2").Value2 = 293.013

# New row 3: "Full-alt" model results (alternate covariates model).
$ws.Range("A3").Value2 = "Full-alt"
$ws.Range("A3").Font.Bold = $true

$ws.Range("D3").Value2 = 70.171
$ws.Range("E3").Value2 = 293.735
$ws.Range("B3").Formula = "=E3-D3"
$ws.Range("C3").Formula = "=B3-D3"

# Move the selection down to A4, matching the post-edit cursor position.
$ws.Range("A4").Select()
